$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Calibration data table: header in row 1, readings in rows 2-8 (A:D).
# Re-sort the readings into chronological order by column A ("time (s)").
$data = @()
for ($r = 2; $r -le 8; $r++) {
    $row = @(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2
    )
    $data += ,$row
}

$sorted = $data | Sort-Object { $_[0] }

for ($i = 0; $i -lt $sorted.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $sorted[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $sorted[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $sorted[$i][2]
    $ws.Cells.Item($r, 4).Value2 = $sorted[$i][3]
}
